$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - rename D1, add new E1 header
$ws.Range("D1").Value = "Lebesgue (int)"
$ws.Range("E1").Value = "Lebesgue (at x0)"

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("B4").Value = 0.1568
$ws.Range("C4").Value = 0.9796
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 1

# Row 5 (B5 already has a value/style in the workbook)
$ws.Range("C5").Value = 0.7053
$ws.Range("D5").Value = 8.7219
$ws.Range("E5").Value = 0.956

# Row 6
$ws.Range("B6").Value = 0.00189775292243155
$ws.Range("C6").Value = 0.5129
$ws.Range("D6").Value = 7.0507
$ws.Range("E6").Value = 0.7218

# Row 7
$ws.Range("B7").Value = 0.000650540124930598
$ws.Range("B7").NumberFormat = "0.00E+00"
$ws.Range("C7").Value = 0.4059
$ws.Range("D7").Value = 5.6471
$ws.Range("E7").Value = 0.490951576804606

# Column D width (closest reachable quantized width to the authored 12.453125)
$ws.Columns.Item(4).ColumnWidth = 11.6

# Final selection, matching the saved cursor position
[void]$ws.Range("B7").Select()
